$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - MAE
$ws.Range("B2").Value = 0.184
$ws.Range("C2").Value = 0.212
$ws.Range("D2").Value = 0.189
$ws.Range("E2").Value = 0.371
$ws.Range("F2").Value = 0.229

# Row 3 - MSE
$ws.Range("B3").Value = 0.112
$ws.Range("C3").Value = 0.111
$ws.Range("D3").Value = 0.108
$ws.Range("E3").Value = 0.626
$ws.Range("F3").Value = 0.157

# Row 5 - mean Y-predicted
$ws.Range("B5").Value = 18.271
$ws.Range("C5").Value = 15.341
$ws.Range("D5").Value = 13.129
$ws.Range("E5").Value = 31.143
$ws.Range("F5").Value = 18.184

# Row 6 - R2
$ws.Range("B6").Value = 0.988
$ws.Range("C6").Value = 0.992
$ws.Range("D6").Value = 0.969
$ws.Range("E6").Value = 0.985
$ws.Range("F6").Value = 0.992
